$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2021 column (S) of data, mirroring column R's styling.
$ws.Range("S4").Value = 2021
$ws.Range("S5").Value = 6.1
$ws.Range("S6").Value = 1.6
$ws.Range("S7").Value = 3.6
$ws.Range("S8").Value = 27.2
$ws.Range("S9").Value = 7.2
$ws.Range("S10").Value = 2.6
$ws.Range("S11").Value = 12.5
$ws.Range("S12").Value = 6.4
$ws.Range("S13").Value = 5.2
$ws.Range("S14").Value = 0.9

# Copy the style from column R onto the new column S cells (rows 4-14).
$ws.Range("R4:R14").Copy()
$ws.Range("S4:S14").PasteSpecial(-4122)  # xlPasteFormats

# Update the active selection to match the post-edit state.
$ws.Range("Q19").Select()
